$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.667.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "156.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,523.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "96.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.14%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.498"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.29%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.659.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.312.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.804"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.589.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.024.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("E45").Value = "  +10.88%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
